# Apply the "Implemented 12/1 feedback, simplified witnesses section" edit.
$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) The first list item ("Print the form(s) ... for more detailed
#    instructions.") gets a brand-new conditional sentence about the
#    child's information prepended in front of it, and the closing
#    sentence of the {% else %} branch changes from "form" to "forms".
#    A hidden "_GoBack" bookmark sits in the middle of this paragraph
#    and must be preserved, so the paragraph is edited in two pieces,
#    split right at the bookmark, instead of one blind whole-text
#    replace.
# ---------------------------------------------------------------------
$goBack = $d.Bookmarks("_GoBack")

# -- piece before the bookmark --------------------------------------
$beforeBm = $d.Range(0, $goBack.Start)
$oldBeforeBm = "{% if other_marriage == False %}Print the form. Read page 2 of the form for more detailed instructions.{% else %}Print the forms. Read page 2 of the forms"
$newBeforeBm = "{% if child_born == False %}Once the child is born, add the child’s information to the {% if other_marriage == True %}forms{% else %} form{% endif %}. {% endif %}"
$beforeBm.Find.Execute($oldBeforeBm, $true, $false, $false, $false, $false, $true, 0, $false, $newBeforeBm, 1) | Out-Null

# -- piece after the bookmark -----------------------------------------
$goBack = $d.Bookmarks("_GoBack")
$afterBm = $d.Range($goBack.End, $d.Content.End)
$oldAfterBm = " for more detailed instructions.{% endif %}"
$newAfterBm = "{% if other_marriage == False %}Print the form. Read page 2 of the form for more detailed instructions.{% else %}Print the forms. Read page 2 of the forms for more detailed instructions.{% endif %}"
$afterBm.Find.Execute($oldAfterBm, $true, $false, $false, $false, $false, $true, 0, $false, $newAfterBm, 1) | Out-Null

# ---------------------------------------------------------------------
# 2) The "Sign the form(s) ... in front of the witnesses you listed."
#    sentence is replaced with the new "Both parents must sign ..."
#    wording.
# ---------------------------------------------------------------------
$oldSign = "Sign the {% if other_marriage == False %}form{% else %}forms{% endif %} with the father {% if other_marriage == True %}and your ex-spouse {% endif %} in front of the witnesses you listed."
$newSign = "Both parents must sign the Voluntary Acknowledgement of Paternity in front of witnesses. {% if other_marriage == True %}The mother and the husband or ex-husband must sign the Denial of Paternity in front of witnesses. {% endif %}Witnesses need to sign, date, and complete the witness information section."
$d.Content.Find.Execute($oldSign, $true, $false, $false, $false, $false, $true, 1, $false, $newSign, 2) | Out-Null
